$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before the current column D ("Terms Typically Offered"),
# shifting it to column G. This makes room for the new
# Corequisites / Concurrent / Recommended columns.
$ws.Range("D1:F1").EntireColumn.Insert()

# New column headers
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Row 2 - GSP 530 (no corequisite text existed, fill with NA)
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"

# Row 3 - GSP 532 (had "Corequisite: GSP 530." inside Prerequisites text)
$ws.Range("C3").Value = "OCOB graduate standing or approval from the Associate Dean of OCOB."
$ws.Range("D3").Value = "GSP 530."
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "F "

# Row 4 - GSP 533
$ws.Range("C4").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D4").Value = "GSP 530."
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "W "

# Row 5 - GSP 535
$ws.Range("C5").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D5").Value = "GSP 530."
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "SU "

# Row 6 - GSP 536
$ws.Range("C6").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D6").Value = "GSP 530."
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "F "

# Row 7 - GSP 537 (no corequisite text existed, fill with NA)
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"

# Row 8 - GSP 538
$ws.Range("C8").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D8").Value = "GSP 530."
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F "

# Row 9 - GSP 539
$ws.Range("C9").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D9").Value = "GSP 530."
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "W "

# Row 10 - GSP 540
$ws.Range("C10").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D10").Value = "GSP 530."
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "SP "

# Row 11 - GSP 541
$ws.Range("C11").Value = "OCOB graduate standing or approval from the Associate Dean."
$ws.Range("D11").Value = "GSP 530."
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "SP "

# Row 12 - GSP 591 (no corequisite text existed, fill with NA)
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
